<#
  feat: add 2022-Q4 data

  Before: workbook has two sheets -> "总计" (summary) and "2022-Q3" (fund
  holdings for 2022-Q3).

  After: a new "2022-Q4" sheet is inserted between them holding fresh
  fund-holding data, the old "2022-Q3" fund sheet is preserved unchanged,
  and the "总计" summary sheet gets an extra row for the new quarter.
#>

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q3" sheet so its original fund
# data survives untouched on its own tab. Worksheet.Copy(Before, After)
# with After = itself drops the copy immediately to its right.
# ---------------------------------------------------------------------
$wsOldQ3 = $wb.Worksheets.Item(2)
$wsOldQ3.Copy($null, $wsOldQ3)

# The sheet that used to be "2022-Q3" becomes "2022-Q4" (new data is
# written into it below). Rename it first so the name "2022-Q3" is free
# for the freshly made duplicate.
$wsOldQ3.Name = "2022-Q4"
$wsQ4 = $wsOldQ3

$wsQ3 = $wb.Worksheets.Item(3)
$wsQ3.Name = "2022-Q3"

$wsSummary = $wb.Worksheets.Item(1)

# The brand-new "2022-Q4" sheet uses the same (default) page margins as
# the "总计" sheet, not the ones inherited from the copied "2022-Q3" fund
# sheet - line them up.
$wsQ4.PageSetup.LeftMargin = $wsSummary.PageSetup.LeftMargin
$wsQ4.PageSetup.RightMargin = $wsSummary.PageSetup.RightMargin
$wsQ4.PageSetup.TopMargin = $wsSummary.PageSetup.TopMargin
$wsQ4.PageSetup.BottomMargin = $wsSummary.PageSetup.BottomMargin
$wsQ4.PageSetup.HeaderMargin = $wsSummary.PageSetup.HeaderMargin
$wsQ4.PageSetup.FooterMargin = $wsSummary.PageSetup.FooterMargin

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - the existing data row now
# describes 2022-Q4, and a new row is appended with the old 2022-Q3
# totals that used to live there. $wsSummary.Range("A2") is never
# overwritten, so it stays a valid style-copy source the whole script.
# ---------------------------------------------------------------------
$wsSummary.Range("B2").Value = "2022-Q4"
$wsSummary.Range("C2").Value = 3
$wsSummary.Range("D2").Value = 0.19

$wsSummary.Range("A2").Copy()
$wsSummary.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$wsSummary.Range("A3").Value = 1
$wsSummary.Range("B3").Value = "2022-Q3"
$wsSummary.Range("C3").Value = 2
$wsSummary.Range("D3").Value = 0.12

# ---------------------------------------------------------------------
# Step 3: replace the "2022-Q4" sheet's fund-holdings table with the new
# quarter's numbers (same header/style layout, now with 3 data rows).
# The header/index-column formatting for this table matches the
# "总计" sheet's header/index style (style index "2"), so pull it from
# there directly instead of from $wsQ4 (whose own header is still the
# old style "1" until we repaint it).
# ---------------------------------------------------------------------
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("C2").Value = "弘毅远方港股通智选领航混合A"
$wsQ4.Range("H2").Value = 4

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("C3").Value = "弘毅远方港股通智选领航混合C"
$wsQ4.Range("H3").Value = 4

$wsSummary.Range("A2").Copy()
$wsQ4.Range("A4").PasteSpecial(-4122)  # xlPasteFormats (index-column style)
$wsQ4.Range("A4").Value = 2
$wsQ4.Range("C4").Value = "恒生前海港股通精选混合"
$wsQ4.Range("H4").Value = 8

# Fund-code / scale / position figures are numeric-looking text in the
# source data (keeps leading zeros like "006537" intact) - mark the
# cells as text *before* assigning so they are not coerced to numbers,
# then drop back to the default (unstyled) look those columns use.
$textCells = @("B2", "D2", "E2", "F2", "G2", `
               "B3", "D3", "E3", "F3", "G3", `
               "B4", "D4", "E4", "F4", "G4")
foreach ($addr in $textCells) {
    $wsQ4.Range($addr).NumberFormat = "@"
}

$wsQ4.Range("B2").Value = "011157"
$wsQ4.Range("D2").Value = "2.84"
$wsQ4.Range("E2").Value = "90.73"
$wsQ4.Range("F2").Value = "4.27"
$wsQ4.Range("G2").Value = "0.1213"

$wsQ4.Range("B3").Value = "011158"
$wsQ4.Range("D3").Value = "0.84"
$wsQ4.Range("E3").Value = "90.73"
$wsQ4.Range("F3").Value = "4.27"
$wsQ4.Range("G3").Value = "0.0359"

$wsQ4.Range("B4").Value = "006537"
$wsQ4.Range("D4").Value = "1.10"
$wsQ4.Range("E4").Value = "91.14"
$wsQ4.Range("F4").Value = "2.99"
$wsQ4.Range("G4").Value = "0.0329"

foreach ($addr in $textCells) {
    $wsQ4.Range($addr).Style = "Normal"
}

# Finally repaint the header row + index column with the correct style
# (now that every value is in place, repainting no longer risks being
# clobbered by a later value assignment).
$wsSummary.Range("B1:D1").Copy()
$wsQ4.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$wsSummary.Range("B1").Copy()
$wsQ4.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats (extra header cols)
$wsSummary.Range("A2").Copy()
$wsQ4.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats

# Keep the originally-active "总计" sheet selected (the various sheet
# copy/rename operations above shift Excel's active-sheet cursor around).
$wsSummary.Activate()
